$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph: "With permission from MegaCrit ... train and test a model."
# ------------------------------------------------------------------

# 1) "...I was able to obtain run data for over 150,000 runs.  For the
#    purposes of efficiency and time I restricted the data to
#    approximately 100,000 runs." -> "...I was able to obtain run data
#    for approximately 188,000 runs."
$d.Content.Find.Execute(
    "over 150,000 runs.  For the purposes of efficiency and time I restricted the data to approximately 100,000 runs.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "approximately 188,000 runs.", 2) | Out-Null

# 2) "This gave me approximately 15,000 runs" -> "This gave me about
#    30,000 runs" (the hidden _GoBack bookmark is relocated into
#    "30,000" further below)
$d.Content.Find.Execute(
    "This gave me approximately 15,000 runs with which",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This gave me about 30,000 runs with which", 2) | Out-Null

# ------------------------------------------------------------------
# Paragraph: "While there is no actually limit ..."
# ------------------------------------------------------------------

# 3) "...an average deck size of 29.7 ***ADD BOX AND WHISKER" ->
#    "...an average deck size of 29.9 cards ***ADD BOX AND WHISKER"
$d.Content.Find.Execute(
    "an average deck size of 29.7 ***ADD BOX AND WHISKER",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "an average deck size of 29.9 cards ***ADD BOX AND WHISKER", 2) | Out-Null

# 4) "...the counts of cards in any given run" (end of paragraph) gains
#    a new trailing sentence about total card features.
$d.Content.Find.Execute(
    "the counts of cards in any given run",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "the counts of cards in any given run.  The number of cards implemented into the game when data was collected was 283; each card has an upgraded version as well, bringing the final count to 566 features.", 2) | Out-Null

# ------------------------------------------------------------------
# Move the hidden "_GoBack" bookmark from the end of the card-count
# paragraph to its new position inside "30,000" (right after "30",
# before ",000").
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete() | Out-Null

$fullText = $d.Content.Text
$marker = "This gave me about 30"
$pos = $fullText.IndexOf($marker)
$bmStart = $pos + $marker.Length
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output $d.Paragraphs(7).Range.Text
Write-Output $d.Paragraphs(8).Range.Text
